$d = $word.ActiveDocument

# --- Edit 1: rewrite the DOC_gw bullet's explanation text -------------------
# Locate the paragraph that starts with "DOC_gw" (currently paragraph 13) and
# replace its contents with the updated explanation, preserving the existing
# "DOC_gw" run/proofErr pair and the paragraph's own pPr/attributes.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("DOC_gw")) {
        $target = $d.Paragraphs($i)
        break
    }
}

$newDocGwXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00EC0CE8" w:rsidRDefault="00EC0CE8" w:rsidP="00057912"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>DOC_gw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">= </w:t></w:r><w:r><w:t xml:space="preserve">personal comm. w/ </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ev</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> &amp; Rudolf Jaffe: DOC </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>gw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> is approximately equal to DOC water column. Use in-lake concentrations for </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> inflow concentration</w:t></w:r></w:p>'

$target.Range.InsertXML($newDocGwXml)

# --- Edit 2: add a new bullet after the "Problem that needs resolving" item -
$problemPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Problem that needs resolving")) {
        $problemPara = $d.Paragraphs($i)
        break
    }
}

$problemPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($problemPara.Index + 1)
$newPara.Range.Text = "OC sediment accumulation rates were provided by Evelyn (from sediment traps)"
